$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected Cypher query text for the "SamplesTab" row (B3) ---
# (removes stray spaces before the commas after Breed / Diagnosis aliases)
$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE demo.breed IN ['Akita'] 
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@

# --- Corrected Cypher query text for the "CasesTab" row (B2) ---
# (removes stray spaces before several commas)
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Akita']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID`,
        coalesce(s.clinical_study_designation, '') AS `Study Code`,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
        coalesce(demo.patient_age_at_enrollment, '') AS Age,
        coalesce(demo.sex, '') AS Sex,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
        coalesce(co.cohort_description, '') AS `Cohort`

'@

# Update SamplesTab row first, then CasesTab row, to match the authored edit order.
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B2").Value = $casesQuery

# Reflect the saved selection/view state: cell B2 selected.
$ws.Activate()
$ws.Range("B2").Select()
